{"js": "// Applies the LOT2043.docx content edits described by the commit diff:\n//  1. Ativa\u00e7\u00e3o date 2018 -> 2022\n//  2. Objetivos paragraph (PT) rewritten\n//  3. Objetivos paragraph (EN, italic) rewritten\n//  4. New docente \"5817181 - Valdeir Arantes\" added on its own line\n//  5. Programa paragraph (PT) rewritten\n//  6. Programa paragraph (EN, italic) rewritten\n//  7. M\u00e9todo value rewritten\n//  8. Crit\u00e9rio value rewritten\n//  9. Norma de recupera\u00e7\u00e3o value rewritten\n\nconst body = context.document.body;\n\n// 1. Ativa\u00e7\u00e3o date\nconst ativacaoResults = body.search(\"Ativa\u00e7\u00e3o: 01/01/2018\", { matchCase: true });\nativacaoResults.load(\"text\");\nawait context.sync();\nif (ativacaoResults.items.length > 0) {\n  ativacaoResults.items[0].insertText(\"Ativa\u00e7\u00e3o: 01/01/2022\", \"Replace\");\n}\nawait context.sync();\n\n// Load paragraphs once, text will be stale after edits below but indices stay stable\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// 2. Objetivos (PT) - paragraph right after the \"Objetivos\" heading\nconst objPtBefore =\n  \"Apresentar ao aluno as caracter\u00edsticas da profiss\u00e3o, os conceitos b\u00e1sicos envolvidos e as \u00e1reas de atua\u00e7\u00e3o do Engenheiro Bioqu\u00edmico, e um panorama do mercado de trabalho de engenharia no Brasil. Mostrar como funciona uma ind\u00fastria de bioprocesso, envolvendo instala\u00e7\u00f5es, tipos de bioprocessos e escala de produ\u00e7\u00e3o.\";\nconst objPtAfter =\n  \"Apresentar aos alunos a Engenharia Bioqu\u00edmica, as caracter\u00edsticas da profiss\u00e3o e orientar quanto as atribui\u00e7\u00f5es e as \u00e1reas de atua\u00e7\u00e3o do Engenheiro Bioqu\u00edmico. Al\u00e9m disso, desenvolver nos alunos uma vis\u00e3o macro dos tipos e etapas de um bioprocesso industrial e, por fim, orientar sobre a atua\u00e7\u00e3o do Engenheiro Bioqu\u00edmico na ind\u00fastria, pesquisa e ensino, e empreendedorismo e inova\u00e7\u00e3o em engenharia.\";\n\n// 3. Objetivos (EN, italic)\nconst objEnBefore =\n  \"Present to the student the characteristics of the profession, the basic concepts involved and the areas of expertise of the Biochemical Engineer, and an overview of the engineering work market in Brazil. Show how a bioprocessing industry works, involving facilities, types of bioprocesses and production scale.\";\nconst objEnAfter =\n  \"To present to the Biochemical Engineering student the characteristics of the profession and to guide in relation to the attributes and the action areas of the biochemical engineering. Besides, to develop in the students a macro view of types and stages of an industrial bioprocess and, finally, to guide about the action of the biochemical engineering on the industry, research and teaching, and entrepreneurship and innovation in engineering.\";\n\n// 5. Programa (PT)\nconst progPtBefore =\n  \"1. Hist\u00f3rico da Engenharia Bioqu\u00edmica: intera\u00e7\u00e3o entre ci\u00eancias biol\u00f3gicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnol\u00f3gicos. 2. Mercado de trabalho da Engenharia do Brasil3. Defini\u00e7\u00f5es e conceitos \u2013 processo enzim\u00e1tico, processo fermentativo gen\u00e9rico, agentes de transforma\u00e7\u00e3o, biorreator, mat\u00e9ria prima, tipos de substratos, convers\u00e3o de substrato em produto, tipos de produtos biotecnol\u00f3gicos, recupera\u00e7\u00e3o de produtos, entre outros.4. \u00c1reas de atua\u00e7\u00e3o do Engenheiro Bioqu\u00edmico5. A Ind\u00fastria de Bioprocessos \u2013 tipos de ind\u00fastrias, equipamentos, instala\u00e7\u00f5es, principais opera\u00e7\u00f5es unit\u00e1rias. 6. Escalas de produ\u00e7\u00e3o \u2013 laborat\u00f3rio, piloto, industrial. 7. Estudo de casos (processos biotecnol\u00f3gicos).8. Visitas supervisionadas \u2013 visitas a laborat\u00f3rios e a ind\u00fastria de bioprocesso.\";\nconst progPtAfter =\n  \"1.Hist\u00f3rico da Engenharia Bioqu\u00edmica: intera\u00e7\u00e3o entre ci\u00eancias biol\u00f3gicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnol\u00f3gicos. 2.Mercado de trabalho da Engenharia do Brasil 3.Atribui\u00e7\u00f5es e \u00e1reas de atua\u00e7\u00e3o do Engenheiro Bioqu\u00edmico 4.Defini\u00e7\u00f5es e conceitos \u2013 processo enzim\u00e1tico, processo fermentativo gen\u00e9rico, agentes de transforma\u00e7\u00e3o, biorreator, mat\u00e9ria prima, tipos de substratos, convers\u00e3o de substrato em produto, tipos de produtos biotecnol\u00f3gicos, recupera\u00e7\u00e3o de produtos, entre outros. 5.A Ind\u00fastria de Bioprocessos \u2013 tipos de ind\u00fastrias, equipamentos, instala\u00e7\u00f5es, principais opera\u00e7\u00f5es unit\u00e1rias. 6.Escalas de produ\u00e7\u00e3o \u2013 laborat\u00f3rio, piloto, industrial. 7.Estudo de casos (processos biotecnol\u00f3gicos). 8.Empreendedorismo e Inova\u00e7\u00e3o em Engenharia.9.Visitas supervisionadas \u2013 visitas a laborat\u00f3rios e a ind\u00fastria de bioprocesso.\";\n\n// 6. Programa (EN, italic)\nconst progEnBefore =\n  \"1. History of Biochemical Engineering: interaction between biological sciences and engineering, multidisciplinarity, peculiarities of biotechnological processes.2. Labor market of Engineering in Brazil3. Definitions and concepts - enzymatic process, generic fermentation process, transformation agents, bioreactor, raw material, types of substrates, substrate conversion into product, types of biotechnological products, product recovery, among others.4. Areas of practice of the Biochemical Engineer5. The Bioprocess Industry - types of industries, equipment, facilities, main unit operations.6. Production scales - laboratory, pilot, industrial.7. Case studies (biotechnological processes).8. Supervised visits - visits to laboratories and the bioprocess industry.\";\nconst progEnAfter =\n  \"1.History of the Biochemical Engineering: interaction between biological sciences and engineering, multidisciplinarity, peculiarities of biotechnological processes.2.Job market of Engineering in Brazil3.Attributes and action areas of biochemical engineering4.Definitions and concepts \u2013 enzymatic process, general fermentative process, transformation agents, bioreactor, raw material, types of substrates, conversion of substrate into product, types of biotechnological products, products recovery, between others.5.The Bioprocesses Industry \u2013 types of industries, equipment, installations, main unit operations6.Production scales \u2013 laboratory, pilot, industrial.7.Studies of cases (biotechnological processes).8.Entrepreneurship and Innovation in Engineering.9.Supervised visitation \u2013 visits to laboratories and bioprocess industry\";\n\n// Map of full-paragraph text replacements keyed by the paragraph's current (before) text.\nconst paraReplacements = [\n  [objPtBefore, objPtAfter],\n  [objEnBefore, objEnAfter],\n  [progPtBefore, progPtAfter],\n  [progEnBefore, progEnAfter],\n];\n\nfor (const item of paragraphs.items) {\n  for (const [before, after] of paraReplacements) {\n    if (item.text === before) {\n      item.insertText(after, \"Replace\");\n      break;\n    }\n  }\n}\nawait context.sync();\n\n// 4. Add new docente line \"5817181 - Valdeir Arantes\" after the existing docente\nconst docenteResults = body.search(\"101761 - Arnaldo M\u00e1rcio Ramalho Prata\", {\n  matchCase: true,\n});\ndocenteResults.load(\"text\");\nawait context.sync();\nif (docenteResults.items.length > 0) {\n  const docenteRange = docenteResults.items[0];\n  docenteRange.insertBreak(\"Line\", \"After\");\n  await context.sync();\n\n  const docenteParagraph = docenteRange.paragraphs.getFirst();\n  const endOfParagraph = docenteParagraph.getRange(\"End\");\n  endOfParagraph.insertText(\"5817181 - Valdeir Arantes\", \"Before\");\n  await context.sync();\n}\n\n// 7. M\u00e9todo value\nconst metodoResults = body.search(\n  \"Provas escritas; participa\u00e7\u00e3o e conte\u00fado de trabalho e semin\u00e1rio;\",\n  { matchCase: true }\n);\nmetodoResults.load(\"text\");\nawait context.sync();\nif (metodoResults.items.length > 0) {\n  metodoResults.items[0].insertText(\n    \"O m\u00e9todo utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as compet\u00eancias t\u00e9cnicas relativas ao tema do projeto, bem como compet\u00eancias transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de lideran\u00e7a e capacidade de comunica\u00e7\u00e3o, dentre outras; exerc\u00edcios individuais realizados no decorrer da disciplina; exerc\u00edcios; din\u00e2micas. Para os projetos, os alunos ser\u00e3o divididos em grupos que desenvolver\u00e3o um projeto durante o semestre relacionado a aplica\u00e7\u00f5es dos conceitos abordados \u00e0 um processo, produto ou servi\u00e7o na \u00e1rea de Engenharia de Bioqu\u00edmica e que relacione com a forma\u00e7\u00e3o acad\u00eamica e atribui\u00e7\u00f5es profissionais do Engenheiro Bioqu\u00edmico.\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// 8. Crit\u00e9rio value\nconst criterioResults = body.search(\n  \"A nota (N) ser\u00e1 composta por ao menos uma prova escrita e trabalhos realizados e apresentados durante o semestre. O peso de cada atividade ser\u00e1 definido segundo crit\u00e9rios do professor.Nota m\u00ednima de aprova\u00e7\u00e3o = 5,0\",\n  { matchCase: true }\n);\ncriterioResults.load(\"text\");\nawait context.sync();\nif (criterioResults.items.length > 0) {\n  criterioResults.items[0].insertText(\n    \"A nota (N) ser\u00e1 individual e ser\u00e1 a m\u00e9dia ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avalia\u00e7\u00e3o dos Pares, Apresenta\u00e7\u00e3o de Trabalhos, dentre outros.\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// 9. Norma de recupera\u00e7\u00e3o value\nconst normaResults = body.search(\n  \"M\u00e9dia Final = (N + Prova Recupera\u00e7\u00e3o)/2Nota Final m\u00ednima para aprova\u00e7\u00e3o= 5,0\",\n  { matchCase: true }\n);\nnormaResults.load(\"text\");\nawait context.sync();\nif (normaResults.items.length > 0) {\n  normaResults.items[0].insertText(\n    \"M\u00e9dia Final = (N + Prova Recupera\u00e7\u00e3o)/2\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n", "ps1": "# Applies the LOT2043.docx content edits described by the commit diff:\n#  1. Ativa\u00e7\u00e3o date 2018 -> 2022\n#  2. Objetivos paragraph (PT) rewritten\n#  3. Objetivos paragraph (EN, italic) rewritten\n#  4. New docente \"5817181 - Valdeir Arantes\" added on its own line\n#  5. Programa paragraph (PT) rewritten\n#  6. Programa paragraph (EN, italic) rewritten\n#  7. M\u00e9todo value rewritten\n#  8. Crit\u00e9rio value rewritten\n#  9. Norma de recupera\u00e7\u00e3o value rewritten\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n# 1. Ativa\u00e7\u00e3o date\nReplace-Text \"Ativa\u00e7\u00e3o: 01/01/2018\" \"Ativa\u00e7\u00e3o: 01/01/2022\"\n\n# 2. Objetivos (PT)\nReplace-Text \"Apresentar ao aluno as caracter\u00edsticas da profiss\u00e3o, os conceitos b\u00e1sicos envolvidos e as \u00e1reas de atua\u00e7\u00e3o do Engenheiro Bioqu\u00edmico, e um panorama do mercado de trabalho de engenharia no Brasil. Mostrar como funciona uma ind\u00fastria de bioprocesso, envolvendo instala\u00e7\u00f5es, tipos de bioprocessos e escala de produ\u00e7\u00e3o.\" `\n             \"Apresentar aos alunos a Engenharia Bioqu\u00edmica, as caracter\u00edsticas da profiss\u00e3o e orientar quanto as atribui\u00e7\u00f5es e as \u00e1reas de atua\u00e7\u00e3o do Engenheiro Bioqu\u00edmico. Al\u00e9m disso, desenvolver nos alunos uma vis\u00e3o macro dos tipos e etapas de um bioprocesso industrial e, por fim, orientar sobre a atua\u00e7\u00e3o do Engenheiro Bioqu\u00edmico na ind\u00fastria, pesquisa e ensino, e empreendedorismo e inova\u00e7\u00e3o em engenharia.\"\n\n# 3. Objetivos (EN, italic)\nReplace-Text \"Present to the student the characteristics of the profession, the basic concepts involved and the areas of expertise of the Biochemical Engineer, and an overview of the engineering work market in Brazil. Show how a bioprocessing industry works, involving facilities, types of bioprocesses and production scale.\" `\n             \"To present to the Biochemical Engineering student the characteristics of the profession and to guide in relation to the attributes and the action areas of the biochemical engineering. Besides, to develop in the students a macro view of types and stages of an industrial bioprocess and, finally, to guide about the action of the biochemical engineering on the industry, research and teaching, and entrepreneurship and innovation in engineering.\"\n\n# 4. Add new docente line \"5817181 - Valdeir Arantes\" after the existing docente\n$find = $d.Content.Find\n$find.Execute(\"101761 - Arnaldo M\u00e1rcio Ramalho Prata\") | Out-Null\nif ($find.Found) {\n    $r = $d.Content\n    $r.Start = $find.Parent.Start\n    $r.End = $find.Parent.End\n    $r.Collapse(0)\n    $r.InsertAfter([char]11 + \"5817181 - Valdeir Arantes\")\n}\n\n# 5. Programa (PT)\nReplace-Text \"1. Hist\u00f3rico da Engenharia Bioqu\u00edmica: intera\u00e7\u00e3o entre ci\u00eancias biol\u00f3gicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnol\u00f3gicos. 2. Mercado de trabalho da Engenharia do Brasil3. Defini\u00e7\u00f5es e conceitos \u2013 processo enzim\u00e1tico, processo fermentativo gen\u00e9rico, agentes de transforma\u00e7\u00e3o, biorreator, mat\u00e9ria prima, tipos de substratos, convers\u00e3o de substrato em produto, tipos de produtos biotecnol\u00f3gicos, recupera\u00e7\u00e3o de produtos, entre outros.4. \u00c1reas de atua\u00e7\u00e3o do Engenheiro Bioqu\u00edmico5. A Ind\u00fastria de Bioprocessos \u2013 tipos de ind\u00fastrias, equipamentos, instala\u00e7\u00f5es, principais opera\u00e7\u00f5es unit\u00e1rias. 6. Escalas de produ\u00e7\u00e3o \u2013 laborat\u00f3rio, piloto, industrial. 7. Estudo de casos (processos biotecnol\u00f3gicos).8. Visitas supervisionadas \u2013 visitas a laborat\u00f3rios e a ind\u00fastria de bioprocesso.\" `\n             \"1.Hist\u00f3rico da Engenharia Bioqu\u00edmica: intera\u00e7\u00e3o entre ci\u00eancias biol\u00f3gicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnol\u00f3gicos. 2.Mercado de trabalho da Engenharia do Brasil 3.Atribui\u00e7\u00f5es e \u00e1reas de atua\u00e7\u00e3o do Engenheiro Bioqu\u00edmico 4.Defini\u00e7\u00f5es e conceitos \u2013 processo enzim\u00e1tico, processo fermentativo gen\u00e9rico, agentes de transforma\u00e7\u00e3o, biorreator, mat\u00e9ria prima, tipos de substratos, convers\u00e3o de substrato em produto, tipos de produtos biotecnol\u00f3gicos, recupera\u00e7\u00e3o de produtos, entre outros. 5.A Ind\u00fastria de Bioprocessos \u2013 tipos de ind\u00fastrias, equipamentos, instala\u00e7\u00f5es, principais opera\u00e7\u00f5es unit\u00e1rias. 6.Escalas de produ\u00e7\u00e3o \u2013 laborat\u00f3rio, piloto, industrial. 7.Estudo de casos (processos biotecnol\u00f3gicos). 8.Empreendedorismo e Inova\u00e7\u00e3o em Engenharia.9.Visitas supervisionadas \u2013 visitas a laborat\u00f3rios e a ind\u00fastria de bioprocesso.\"\n\n# 6. Programa (EN, italic)\nReplace-Text \"1. History of Biochemical Engineering: interaction between biological sciences and engineering, multidisciplinarity, peculiarities of biotechnological processes.2. Labor market of Engineering in Brazil3. Definitions and concepts - enzymatic process, generic fermentation process, transformation agents, bioreactor, raw material, types of substrates, substrate conversion into product, types of biotechnological products, product recovery, among others.4. Areas of practice of the Biochemical Engineer5. The Bioprocess Industry - types of industries, equipment, facilities, main unit operations.6. Production scales - laboratory, pilot, industrial.7. Case studies (biotechnological processes).8. Supervised visits - visits to laboratories and the bioprocess industry.\" `\n             \"1.History of the Biochemical Engineering: interaction between biological sciences and engineering, multidisciplinarity, peculiarities of biotechnological processes.2.Job market of Engineering in Brazil3.Attributes and action areas of biochemical engineering4.Definitions and concepts \u2013 enzymatic process, general fermentative process, transformation agents, bioreactor, raw material, types of substrates, conversion of substrate into product, types of biotechnological products, products recovery, between others.5.The Bioprocesses Industry \u2013 types of industries, equipment, installations, main unit operations6.Production scales \u2013 laboratory, pilot, industrial.7.Studies of cases (biotechnological processes).8.Entrepreneurship and Innovation in Engineering.9.Supervised visitation \u2013 visits to laboratories and bioprocess industry\"\n\n# 7. M\u00e9todo value\nReplace-Text \"Provas escritas; participa\u00e7\u00e3o e conte\u00fado de trabalho e semin\u00e1rio;\" `\n             \"O m\u00e9todo utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as compet\u00eancias t\u00e9cnicas relativas ao tema do projeto, bem como compet\u00eancias transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de lideran\u00e7a e capacidade de comunica\u00e7\u00e3o, dentre outras; exerc\u00edcios individuais realizados no decorrer da disciplina; exerc\u00edcios; din\u00e2micas. Para os projetos, os alunos ser\u00e3o divididos em grupos que desenvolver\u00e3o um projeto durante o semestre relacionado a aplica\u00e7\u00f5es dos conceitos abordados \u00e0 um processo, produto ou servi\u00e7o na \u00e1rea de Engenharia de Bioqu\u00edmica e que relacione com a forma\u00e7\u00e3o acad\u00eamica e atribui\u00e7\u00f5es profissionais do Engenheiro Bioqu\u00edmico.\"\n\n# 8. Crit\u00e9rio value\nReplace-Text \"A nota (N) ser\u00e1 composta por ao menos uma prova escrita e trabalhos realizados e apresentados durante o semestre. O peso de cada atividade ser\u00e1 definido segundo crit\u00e9rios do professor.Nota m\u00ednima de aprova\u00e7\u00e3o = 5,0\" `\n             \"A nota (N) ser\u00e1 individual e ser\u00e1 a m\u00e9dia ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avalia\u00e7\u00e3o dos Pares, Apresenta\u00e7\u00e3o de Trabalhos, dentre outros.\"\n\n# 9. Norma de recupera\u00e7\u00e3o value\nReplace-Text \"M\u00e9dia Final = (N + Prova Recupera\u00e7\u00e3o)/2Nota Final m\u00ednima para aprova\u00e7\u00e3o= 5,0\" `\n             \"M\u00e9dia Final = (N + Prova Recupera\u00e7\u00e3o)/2\"\n"}
